$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2075
    $ws.Range("F4").Value = 130
    $ws.Range("F7").Value = 1716
    $ws.Range("F8").Value = 27
    $ws.Range("F9").Value = 708
    $ws.Range("F13").Value = 104
    $ws.Range("F15").Value = 11
    $ws.Range("F18").Value = 137
    $ws.Range("F19").Value = 3977
    $ws.Range("F23").Value = 381
    $ws.Range("F25").Value = 730
    $ws.Range("F27").Value = 12
    $ws.Range("F29").Value = 1778
    $ws.Range("F31").Value = 33
    $ws.Range("F32").Value = 65
}

Write-Host "Done updating F-column values on sheets: $sheetNames"
